# aggiornamento fino a 02/05
# Adds rows 239-244 to Sheet1, continuing the daily series in columns A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style template range (last existing data row, columns A:D only) - copy its
# formats (incl. date number format + borders) down onto the new rows before
# writing values. Restricting to A:D avoids touching the rest of the row.
$templateRange = $ws.Range("A238:D238")

$newData = @(
    @(239, 44313, 0, 3, 187.0324189526185),
    @(240, 44314, 0, 3, 187.0324189526185),
    @(241, 44315, 0, 0, 0),
    @(242, 44316, 0, 0, 0),
    @(243, 44317, 0, 0, 0),
    @(244, 44318, 1, 1, 62.34413965087282)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $templateRange.Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}

$excel.CutCopyMode = 0
